$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 10) of arrival data, mirroring the structure of the
# existing rows (flight W92182 from London, same as row 3, but with a new
# status/difference reflecting an additional arrival on Sunday, Jan 15).
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Sunday, Jan 15"
$ws.Range("C10").Value = "3:40 PM"
$ws.Range("D10").Value = "W92182"
$ws.Range("E10").Value = "London"
$ws.Range("F10").Value = "(LTN)"
$ws.Range("G10").Value = "Wizz Air "
$ws.Range("H10").Value = "A320"
$ws.Range("I10").Value = "(G-WUKF)"
$ws.Range("J10").Value = "3:06 PM"
$ws.Range("L10").Value = "0 hours, -34 minutes"
